$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 708 (shifts existing rows 708+ down by one)
$ws.Rows.Item(708).Insert()

# Populate the newly inserted row 708 with the new data point.
# Force column A to Text format first so the date-like string is not
# auto-converted into a date serial value, then clear the formatting back
# to the sheet default (matches the unstyled neighboring data cells).
$ws.Cells.Item(708, 1).NumberFormat = "@"
$ws.Cells.Item(708, 1).Value = "2026/01/23"
$ws.Cells.Item(708, 1).ClearFormats()
$ws.Cells.Item(708, 2).Value = "金"
$ws.Cells.Item(708, 3).Value = 7
$ws.Cells.Item(708, 4).Value = 178
